$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 1).Value = "平潭发展"
$ws.Cells.Item(2, 2).Value = "平潭发展"
$ws.Cells.Item(2, 3).Value = "襄阳轴承"

$ws.Cells.Item(3, 1).Value = "航天发展"
$ws.Cells.Item(3, 2).Value = "航天发展"
$ws.Cells.Item(3, 3).Value = "平潭发展"

$ws.Cells.Item(4, 1).Value = "乾照光电"
$ws.Cells.Item(4, 2).Value = "雷科防务"
$ws.Cells.Item(4, 3).Value = "航天发展"

$ws.Cells.Item(5, 1).Value = "雷科防务"
$ws.Cells.Item(5, 2).Value = "乾照光电"
$ws.Cells.Item(5, 3).Value = "国晟科技"

$ws.Cells.Item(6, 1).Value = "长盈精密"
$ws.Cells.Item(6, 2).Value = "实达集团"
$ws.Cells.Item(6, 3).Value = "实达集团"

$ws.Cells.Item(7, 1).Value = "国晟科技"
$ws.Cells.Item(7, 2).Value = "海欣食品"
$ws.Cells.Item(7, 3).Value = "天风证券"

$ws.Cells.Item(8, 1).Value = "实达集团"
$ws.Cells.Item(8, 2).Value = "国晟科技"
$ws.Cells.Item(8, 3).Value = "合富中国"

$ws.Cells.Item(9, 1).Value = "襄阳轴承"
$ws.Cells.Item(9, 2).Value = "襄阳轴承"
$ws.Cells.Item(9, 3).Value = "雷科防务"

$ws.Cells.Item(10, 1).Value = "通宇通讯"
$ws.Cells.Item(10, 2).Value = "福日电子"
$ws.Cells.Item(10, 3).Value = "海欣食品"

$ws.Cells.Item(11, 1).Value = "海欣食品"
$ws.Cells.Item(11, 2).Value = "安妮股份"
$ws.Cells.Item(11, 3).Value = "乾照光电"

$ws.Cells.Item(12, 1).Value = "合富中国"
$ws.Cells.Item(12, 2).Value = "合富中国"
$ws.Cells.Item(12, 3).Value = "万通发展"

$ws.Cells.Item(13, 1).Value = "海峡创新"
$ws.Cells.Item(13, 2).Value = "通宇通讯"
$ws.Cells.Item(13, 3).Value = "海峡创新"

$ws.Cells.Item(14, 1).Value = "福日电子"
$ws.Cells.Item(14, 2).Value = "长盈精密"
$ws.Cells.Item(14, 3).Value = "新华都"

$ws.Cells.Item(15, 1).Value = "多氟多"
$ws.Cells.Item(15, 2).Value = "海峡创新"
$ws.Cells.Item(15, 3).Value = "航天动力"

$ws.Cells.Item(16, 1).Value = "航天动力"
$ws.Cells.Item(16, 2).Value = "多氟多"
$ws.Cells.Item(16, 3).Value = "特发信息"

$ws.Cells.Item(17, 1).Value = "赛微电子"
$ws.Cells.Item(17, 2).Value = "广汽集团"
$ws.Cells.Item(17, 3).Value = "多氟多"

$ws.Cells.Item(18, 1).Value = "蓝色光标"
$ws.Cells.Item(18, 2).Value = "蓝色光标"
$ws.Cells.Item(18, 3).Value = "海王生物"

$ws.Cells.Item(19, 1).Value = "茂业商业"
$ws.Cells.Item(19, 2).Value = "海南瑞泽"
$ws.Cells.Item(19, 3).Value = "梅雁吉祥"

$ws.Cells.Item(20, 1).Value = "安妮股份"
$ws.Cells.Item(20, 2).Value = "航天动力"
$ws.Cells.Item(20, 3).Value = "长盈精密"

$ws.Cells.Item(21, 1).Value = "顺灏股份"
$ws.Cells.Item(21, 2).Value = "特一药业"
$ws.Cells.Item(21, 3).Value = "众生药业"
